# Fill in the "mode" and "median" sheets with the computed summary-statistic
# rows (rows 2-19) for the three groups (Unmanipulated / NoMAP / MAP) that
# already exist, fully populated, on the "sd", "range", "iqr" and "mean"
# sheets. Mirrors the pattern already present on those sibling sheets,
# including the columns whose series end early (col A stops after a few
# rows, col B a little later, col C runs the full 19 rows).

$wb = $excel.ActiveWorkbook

# row, colA, colB, colC ($null = leave blank, series ended)
$mode_data = @(
  @(2, 4, 6, 33),
  @(3, 4, 7, 3),
  @(4, 4, 9, 5),
  @(5, 3, 9, 4),
  @(6, 4, 4, 49),
  @(7, 5, 9, 5),
  @(8, $null, 15, 4),
  @(9, $null, $null, 4),
  @(10, $null, $null, 13),
  @(11, $null, $null, 6),
  @(12, $null, $null, 7),
  @(13, $null, $null, 3),
  @(14, $null, $null, 17),
  @(15, $null, $null, 4),
  @(16, $null, $null, 15),
  @(17, $null, $null, 9),
  @(18, $null, $null, 20),
  @(19, $null, $null, 8)
)

$median_data = @(
  @(2, 4.62, 7.26, 34.76),
  @(3, 5.72, 6.16, 3.08),
  @(4, 4.4, 10.56, 5.5),
  @(5, 4.84, 8.36, 6.38),
  @(6, 6.16, 4.4, 48.84),
  @(7, $null, 13.64, 5.72),
  @(8, $null, 12.32, 3.96),
  @(9, $null, $null, 4.62),
  @(10, $null, $null, 18.92),
  @(11, $null, $null, 6.82),
  @(12, $null, $null, 9.46),
  @(13, $null, $null, 3.08),
  @(14, $null, $null, 15.84),
  @(15, $null, $null, 4.4),
  @(16, $null, $null, 18.26),
  @(17, $null, $null, 11.44),
  @(18, $null, $null, 36.96),
  @(19, $null, $null, 9.02)
)

function Fill-Sheet($sheetName, $rows) {
  $ws = $wb.Worksheets.Item($sheetName)
  foreach ($entry in $rows) {
    $r = $entry[0]
    if ($null -ne $entry[1]) { $ws.Cells.Item($r, 1).Value = $entry[1] }
    if ($null -ne $entry[2]) { $ws.Cells.Item($r, 2).Value = $entry[2] }
    if ($null -ne $entry[3]) { $ws.Cells.Item($r, 3).Value = $entry[3] }
  }
}

Fill-Sheet "mode" $mode_data
Fill-Sheet "median" $median_data
